$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Fix typo in heading: "Description" -> "description" (this broke XLS file chunking)
$ws.Range("G3").Value = "description"

# Move the active cell selection from G5 to G4
$ws.Range("G4").Select()
